{"js": "// Replace hardcoded date / of\u00edcio number / ano with the report's\n// placeholder fields (dataEnvio, numeroOficio, anoParecer) and spell\n// out \"CONSEPE\", per commit: \"fix: colocado dataEnvio e numeroOficio no docx\".\n\nconst body = context.document.body;\n\n// 1) \"... no dia 18 de novembro de 2024, por meio do OF\u00cdCIO CIRCULAR N\u00ba 2/2024/SE/CAMTUC/UFPA, sendo identificada ...\"\n//    -> \"... no dia dataEnvio, por meio do OF\u00cdCIO CIRCULAR N\u00ba numeroOficio/anoParecer/SE/CAMTUC/UFPA, sendo identificada ...\"\nconst firstMention = body.search(\n  \"18 de novembro de 2024, por meio do OF\u00cdCIO CIRCULAR N\u00ba 2/2024/SE/CAMTUC/UFPA\",\n  { matchCase: true }\n);\nfirstMention.load(\"items\");\n\n// 2) \"Constam como documentos da proposta: (a) OF\u00cdCIO CIRCULAR N\u00ba 2/2024/SE/CAMTUC/UFPA;\"\n//    -> \"Constam como documentos da proposta: (a) OF\u00cdCIO CIRCULAR N\u00ba numeroOficio/anoParecer/SE/CAMTUC/UFPA;\"\nconst secondMention = body.search(\n  \"Constam como documentos da proposta: (a) OF\u00cdCIO CIRCULAR N\u00ba 2/2024/SE/CAMTUC/UFPA;\",\n  { matchCase: true }\n);\nsecondMention.load(\"items\");\n\n// 3) \"... Pesquisa e Extens\u00e3o (CONSEP), resolu\u00e7\u00e3o ...\" -> \"... (CONSEPE), resolu\u00e7\u00e3o ...\"\nconst consep = body.search(\"Extens\u00e3o (CONSEP), resolu\u00e7\u00e3o\", { matchCase: true });\nconsep.load(\"items\");\n\nawait context.sync();\n\nif (firstMention.items.length > 0) {\n  firstMention.items[0].insertText(\n    \"18 de novembro de 2024, por meio do OF\u00cdCIO CIRCULAR N\u00ba 2/2024/SE/CAMTUC/UFPA\"\n      .replace(\"18 de novembro de 2024\", \"dataEnvio\")\n      .replace(\"2/2024\", \"numeroOficio/anoParecer\"),\n    Word.InsertLocation.replace\n  );\n}\n\nif (secondMention.items.length > 0) {\n  secondMention.items[0].insertText(\n    \"Constam como documentos da proposta: (a) OF\u00cdCIO CIRCULAR N\u00ba 2/2024/SE/CAMTUC/UFPA;\".replace(\n      \"2/2024\",\n      \"numeroOficio/anoParecer\"\n    ),\n    Word.InsertLocation.replace\n  );\n}\n\nif (consep.items.length > 0) {\n  consep.items[0].insertText(\"Extens\u00e3o (CONSEPE), resolu\u00e7\u00e3o\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace hardcoded date / of\u00edcio number / ano with the report's\n# placeholder fields (dataEnvio, numeroOficio, anoParecer) and spell\n# out \"CONSEPE\", per commit: \"fix: colocado dataEnvio e numeroOficio no docx\".\n\n$d = $word.ActiveDocument\n\n# 1) \"... no dia 18 de novembro de 2024, por meio do OF\u00cdCIO CIRCULAR N\u00ba 2/2024/SE/CAMTUC/UFPA, sendo identificada ...\"\n#    -> \"... no dia dataEnvio, por meio do OF\u00cdCIO CIRCULAR N\u00ba numeroOficio/anoParecer/SE/CAMTUC/UFPA, sendo identificada ...\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"18 de novembro de 2024, por meio do OF\u00cdCIO CIRCULAR N\u00ba 2/2024/SE/CAMTUC/UFPA\"\n$find.Replacement.Text = \"dataEnvio, por meio do OF\u00cdCIO CIRCULAR N\u00ba numeroOficio/anoParecer/SE/CAMTUC/UFPA\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 2) \"Constam como documentos da proposta: (a) OF\u00cdCIO CIRCULAR N\u00ba 2/2024/SE/CAMTUC/UFPA;\"\n#    -> \"Constam como documentos da proposta: (a) OF\u00cdCIO CIRCULAR N\u00ba numeroOficio/anoParecer/SE/CAMTUC/UFPA;\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Constam como documentos da proposta: (a) OF\u00cdCIO CIRCULAR N\u00ba 2/2024/SE/CAMTUC/UFPA;\"\n$find2.Replacement.Text = \"Constam como documentos da proposta: (a) OF\u00cdCIO CIRCULAR N\u00ba numeroOficio/anoParecer/SE/CAMTUC/UFPA;\"\n$find2.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 3) \"... Pesquisa e Extens\u00e3o (CONSEP), resolu\u00e7\u00e3o ...\" -> \"... (CONSEPE), resolu\u00e7\u00e3o ...\"\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Text = \"Extens\u00e3o (CONSEP), resolu\u00e7\u00e3o\"\n$find3.Replacement.Text = \"Extens\u00e3o (CONSEPE), resolu\u00e7\u00e3o\"\n$find3.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n"}
